$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = 41949
$ws.Range("B106").Value = 0.83611111111111114
$ws.Range("C106").Value = 0.92013888888888884
$ws.Range("D106").Value = 5
$ws.Range("E106").Formula = "=IF(AND(NOT(ISBLANK(B106)),NOT(ISBLANK(C106))), (C106-B106) * 24 - D106/60, """")"
$ws.Range("F106").Value = "Coding"

$ws.Range("A107").Select()
